# WorkCenters.xlsx fixes
#  - Rename two BAXTER SPRINGS work-center descriptions (C60, C61)
#  - Append two new PAWNEE CITY work-center rows (71, 72)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename existing descriptions ---
$ws.Range("C60").Value = "BAXTER BIG BAG - LINE 2"
$ws.Range("C61").Value = "BAXTER SMALL BAG -  LINE 1"

# --- New row 71 ---
# Facility/WorkCenter codes are stored as text in this sheet (leading
# apostrophe keeps "20005"/"74510" as text instead of being coerced to a
# number), matching the existing A:B columns.
$ws.Range("A71").Value = "'20005"
$ws.Range("B71").Value = "'74510"
$ws.Range("C71").Value = "SUP PAWNEE CITY"
$ws.Range("D71").Value = "PACK"

# --- New row 72 ---
$ws.Range("A72").Value = "'20005"
$ws.Range("B72").Value = "'74511"
$ws.Range("C72").Value = "PSGLEE PAWNEE CITY"
$ws.Range("D72").Value = "PACK"
